$wb = $excel.ActiveWorkbook

# The "House part sizes" sheet/table is being retired - its single
# "Square meters" measurement per part is folded into the
# "House parts and materials" table as a new "Size" column instead, so the
# whole material list lives in one populated table.
$wsSizes = $wb.Worksheets.Item("House part sizes")

# Remember each part's size (by name, so row order doesn't matter) before
# the source sheet disappears.
$sizeByPart = @{}
for ($r = 2; $r -le 4; $r++) {
    $partName = $wsSizes.Cells.Item($r, 1).Value()
    $sizeValue = $wsSizes.Cells.Item($r, 2).Value()
    $sizeByPart[$partName] = $sizeValue
}

[void]$wsSizes.Delete()

# Re-resolve the materials sheet now that the workbook only has two tabs.
$wsMaterials = $wb.Worksheets.Item("House parts and materials")

# Grow the MaterialList table by one column to hold the new "Size" data.
$lo = $wsMaterials.ListObjects.Item("MaterialList")
[void]$lo.ListColumns.Add()
$wsMaterials.Range("E1").Value = "Size"

for ($r = 2; $r -le 4; $r++) {
    $partName = $wsMaterials.Cells.Item($r, 1).Value()
    $wsMaterials.Cells.Item($r, 5).Value = $sizeByPart[$partName]
}

# Match the formatting already used by the rest of the table (bold/centered
# header, vertically-centered + wrapped data cells).
[void]$wsMaterials.Range("D1").Copy()
[void]$wsMaterials.Range("E1").PasteSpecial(-4122)
[void]$wsMaterials.Range("D2:D4").Copy()
[void]$wsMaterials.Range("E2:E4").PasteSpecial(-4122)

# The paste brings over D's text/value too, so re-apply the header label and
# the numeric sizes after the formatting-only paste.
$wsMaterials.Range("E1").Value = "Size"
for ($r = 2; $r -le 4; $r++) {
    $partName = $wsMaterials.Cells.Item($r, 1).Value()
    $wsMaterials.Cells.Item($r, 5).Value = $sizeByPart[$partName]
}

[void]$wsMaterials.Activate()
[void]$wsMaterials.Range("E4").Select()
